$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 5 values
$ws.Range("B5").Value = "5001404-20.2017.8.21.0042"
$ws.Range("C5").Value = "0059732-52.2019.8.21.9000"
$ws.Range("D5").Value = "CIV.12711.02"

# Update row 6 values
$ws.Range("B6").Value = "5001392-06.2017.8.21.0042"
$ws.Range("C6").Value = "0038543-18.2019.8.21.9000"
$ws.Range("D6").Value = "CIV.14284.02"

# Delete rows 7 through 12 (shift cells up, xlShiftUp = -4162)
$ws.Range("A7:E12").Delete(-4162)
